$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit text format
# applied first (then reset back to Normal style) so Excel keeps them as text
# instead of silently converting them to numeric values.
$textCells = @('D5', 'D6', 'D10', 'D11', 'D13', 'D14', 'D19', 'D20', 'D21', 'D23', 'D24', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '52.438.32'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '2.924.27'
$ws.Range('E3').Value = '  +4.73%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '352.77'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').Value = '112.86'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').Value = '40.37'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.0863'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = '20.22'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '7.86'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '3.383.79'
$ws.Range('E15').Value = '  +4.58%  '
$ws.Range('D16').Value = '2.937.30'
$ws.Range('E16').Value = '  +4.78%  '
$ws.Range('E17').Value = '  +6.69%  '
$ws.Range('D18').Value = '52.425.58'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').Value = '7.77'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '3.37'
$ws.Range('E20').Value = '  +6.61%  '
$ws.Range('D21').Value = '14.47'
$ws.Range('E21').Value = '  +8.09%  '
$ws.Range('D22').Value = '0.0₃0984'
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').Value = '71.26'
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('D24').Value = '271.82'
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('D26').Value = '27.00'
$ws.Range('E26').Value = '  +3.63%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '10.59'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('D30').Value = '38.58'
$ws.Range('E30').Value = '  +3.88%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '2.25'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '6.47'
$ws.Range('E32').Value = '  +3.36%  '
$ws.Range('D33').Value = '6.18'
$ws.Range('E33').Value = '  +8.73%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '53.08'
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0940'
$ws.Range('E35').Value = '  +10.36%  '
$ws.Range('D36').Value = '0.0457'
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '3.37'
$ws.Range('E38').Value = '  +7.67%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '2.09'
$ws.Range('E39').Value = '  +5.82%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '18.89'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  +7.56%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '24.12'
$ws.Range('E42').Value = '  +10.57%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  +2.26%  '
$ws.Range('D44').Value = '121.93'
$ws.Range('E44').Value = '  +2.33%  '
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').Value = '3.58'
$ws.Range('E46').Value = '  +5.19%  '
$ws.Range('D47').Value = '2.213.02'
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('D48').Value = '2.51'
$ws.Range('E48').Value = '  +6.64%  '
$ws.Range('D49').Value = '0.267'
$ws.Range('E49').Value = '  +23.12%  '
$ws.Range('D50').Value = '0.957'
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').Value = '0.0327'
$ws.Range('E51').Value = '  +13.28%  '

# Reset style back to Normal for cells we forced to text format, so no
# stray style index is introduced relative to the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}